# "Historia de Usuario Épica.docx" - commit: "Clase 27/08/2024 avance
# creacion base de datos"
#
# Body-text clean-up: several paragraphs had their text split across
# multiple runs (left over from copy/paste + manual retyping with no
# formatting differences); this collapses each back down to a single run.
# It also fixes "Consumidor" -> "Cliente" in the data-model paragraph and
# removes the unfinished "13..16" placeholder section at the end of the
# document.

$d = $word.ActiveDocument

$wdFindContinue  = 1
$wdReplaceAll    = 2

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, `
                             $true, $wdFindContinue, $false, $replace, $wdReplaceAll) | Out-Null
}

# --- Re-type paragraphs whose text is unchanged but was spread across
#     several runs, so they collapse back into a single run each. ---

Replace-All "Historia de Usuario Épica" "Historia de Usuario Épica"
Replace-All "1. Nombre del proyecto: Envía" "1. Nombre del proyecto: Envía"
Replace-All "2. Proceso: Envió de paquetes" "2. Proceso: Envió de paquetes"

$mision = "5. Misión: La en empresa envía se encarga de la entrega de paquetes de manera eficiente a diferente cliente de manera regional."
Replace-All $mision $mision

$vision = "6.Vision: La empresa envía se proyecta como una empresa reconocida a nivel nacional por su calidad en el servicio."
Replace-All $vision $vision

# --- Data-model paragraph: "Consumidor" -> "Cliente" ---
Replace-All "Consumidor" "Cliente"

# --- Drop the trailing, never-written "13..16" section. It starts right
#     after "12. Pasos del proceso:" last bullet and runs to the end of
#     the document, so just cut from its first heading through the end of
#     the story. ---
$startMarker = $d.Content
$found = $startMarker.Find.Execute("13. Objetos de Alto Valor", $false, $false, $false, `
                                    $false, $false, $true, $wdFindContinue, $false, "", 0)
if ($found) {
    $deleteStart = $startMarker.Start
    $deleteEnd = $d.Content.End
    $tailRange = $d.Range($deleteStart, $deleteEnd)
    $tailRange.Delete() | Out-Null
}
